$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text (German gender-inclusive wording) and make all
# text headers bold, matching the updated "Fabrikam Q1 marketing
# campaigns" header row / Table1 column names.
$ws.Range("A1").Value = "Kampagnenbesitzer/Kampagnenbesitzerin"
$ws.Range("G1").Value = "Gesamtanzahl der Zielbenutzenden"
$ws.Range("H1").Value = "Aktive Benutzende"

$ws.Range("A1:H1").Font.Bold = $true
$ws.Range("A1:H1").Font.Color = 16777215

# Keep the Table1 column headers in sync with the worksheet header row.
$table = $ws.ListObjects.Item("Table1")
$table.ListColumns.Item(1).Name = "Kampagnenbesitzer/Kampagnenbesitzerin"
$table.ListColumns.Item(7).Name = "Gesamtanzahl der Zielbenutzenden"
$table.ListColumns.Item(8).Name = "Aktive Benutzende"
